$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: remove "Version" from column A and shift the remaining
# header cells (Code, Description, Validation script) one column left.
$ws.Range("A2").Value = $ws.Range("B2").Value2
$ws.Range("B2").Value = $ws.Range("C2").Value2
$ws.Range("C2").Value = $ws.Range("D2").Value2
$ws.Range("D2").Clear()

# Row 3: remove "1" (the version number value) from column A and
# shift the remaining cells (ATTACHMENT, Attachment) one column left.
$ws.Range("A3").Value = $ws.Range("B3").Value2
$ws.Range("B3").Value = $ws.Range("C3").Value2
$ws.Range("C3").Clear()

# Row 3 now behaves like a normal data row (same height as the others).
$ws.Rows.Item(3).RowHeight = 15

# Move the active selection to A2 (previously A4).
[void]$ws.Range("A2").Select()
